$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values per row, per repulled data / mean calculation fix.
$updates = @{
    3  = -4
    6  = -4
    7  = 8
    8  = 1
    9  = -1
    10 = -3
    12 = 1
    13 = -1
    14 = 1
    15 = 2
    16 = 7
    17 = -3
    18 = -2
    19 = -2
    22 = -6
    23 = -2
    24 = -1
    25 = -2
    26 = 1
    27 = -4
    28 = 2
    29 = -4
    30 = 1
    31 = -3
    32 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
